$d = $word.ActiveDocument

# wdAlignParagraphCenter = 1, wdAlignParagraphLeft = 0
$wdAlignParagraphLeft = 0
$wdAlignParagraphCenter = 1

# The document currently has a single empty paragraph.
# Turn it into the centered "Capítulo 1." heading line, then add two
# more centered empty paragraphs, and finally a left-aligned paragraph
# that just contains a tab character.

$p1 = $d.Paragraphs.Item(1)
$p1.Range.Text = "Capítulo 1."
$p1.Alignment = $wdAlignParagraphCenter

$r = $d.Paragraphs.Item(1).Range
$r.Collapse(0)
$r.InsertParagraphAfter()

$p2 = $d.Paragraphs.Item(2)
$p2.Alignment = $wdAlignParagraphCenter

$r2 = $p2.Range
$r2.Collapse(0)
$r2.InsertParagraphAfter()

$p3 = $d.Paragraphs.Item(3)
$p3.Alignment = $wdAlignParagraphCenter

$r3 = $p3.Range
$r3.Collapse(0)
$r3.InsertParagraphAfter()

$p4 = $d.Paragraphs.Item(4)
$p4.Alignment = $wdAlignParagraphLeft
$p4.Range.Text = "`t"
